# orders.xlsx: the order log effectively drops its oldest entry (1/6/2020,
# East, Jones, Pencil) and its newest entry (12/21/2021, Central, Andrews,
# Binder) was removed as well, so every remaining record shifts up by one
# row. Along the way a data-entry correction is made: the "Parent / Pen"
# order date is fixed from 11/8/2020 to 11/8/2021.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (old row 44) first, while row indices still
# refer to the original layout.
$ws.Rows(44).Delete()

# Remove the first data row (old row 2). Everything below shifts up by
# one row, matching the target layout (A1:G42).
$ws.Rows(2).Delete()

# Correct the order date that lands on the new row 19 (was row 20):
# 11/8/2020 -> 11/8/2021. Use a leading apostrophe so Excel keeps the
# value as literal text (matching every other date in the sheet) instead
# of re-interpreting it as a date serial number.
$ws.Range("A19").Value = "'11/8/2021"

Write-Host "done"
